$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; temporarily unprotect so the cell values can be
# updated, then re-protect it afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidentiality / disclosure note (A16).
$ws.Range("A16").Value2 = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."
# Setting a two-line value auto-expands the row height; restore it to fit
# the (unchanged) default row height so the row stays at its original size.
$ws.Rows.Item(16).AutoFit()

# Refresh the Weight (D) and Percent Change (E) columns for rows 2-13.
$ws.Range("D2").Value2 = 0.03066642390291778
$ws.Range("E2").Value2 = -0.005066956207021356

$ws.Range("D3").Value2 = 0.02418937068902089
$ws.Range("E3").Value2 = 0.005330748727889523

$ws.Range("D4").Value2 = 0.05309392837056059
$ws.Range("E4").Value2 = -0.003593675131768115

$ws.Range("D5").Value2 = 0.1378356032159546
$ws.Range("E5").Value2 = 0.008207934336525335

$ws.Range("D6").Value2 = 0.03110925989896206
$ws.Range("E6").Value2 = -0.006746626686656598

$ws.Range("D7").Value2 = 0.1200191550574465
$ws.Range("E7").Value2 = 0.007389812615465896

$ws.Range("D8").Value2 = 0.1025975305439612
$ws.Range("E8").Value2 = 0.00019204916458615

$ws.Range("D9").Value2 = 0.0282260246895721
$ws.Range("E9").Value2 = -0.002929247408742675

$ws.Range("D10").Value2 = 0.1229365600958751
$ws.Range("E10").Value2 = -0.006897950902820082

$ws.Range("D11").Value2 = 0.2467298600658206
$ws.Range("E11").Value2 = 0.01661384807452237

$ws.Range("D12").Value2 = 0.1025962834699087
$ws.Range("E12").Value2 = 0.007147232999795916

$ws.Range("D13").Value2 = 1
$ws.Range("E13").Value2 = 0.005512566141690378

# Restore sheet protection.
$ws.Protect()
